# Auto commit update: refresh Metrics values (B2:B13) and let the
# dependent "today" sheet formulas (B11:B22, E11:E22, F11:F22) and the
# TODAY()-1 cell (A1) recalculate automatically.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# New metric values (Metrics!B2:B13)
$wsMetrics.Range("B2").Value  = 377952.90000000008
$wsMetrics.Range("B3").Value  = 332945.24000000005
$wsMetrics.Range("B4").Value  = 116401.36
$wsMetrics.Range("B5").Value  = 15396
$wsMetrics.Range("B6").Value  = 5174198.6500000004
$wsMetrics.Range("B7").Value  = 4375021.9200000009
$wsMetrics.Range("B8").Value  = 1523361.1900000004
$wsMetrics.Range("B9").Value  = 201603
$wsMetrics.Range("B10").Value = 33639579.640000015
$wsMetrics.Range("B11").Value = 31650297.080000006
$wsMetrics.Range("B12").Value = 11805083.229999999
$wsMetrics.Range("B13").Value = 1299233

# Recalculate so the "today" sheet's formulas (which reference the
# Metrics sheet) and the TODAY()-1 cell pick up the new values.
$excel.CalculateFull()

# Restore the selections recorded in each sheet's view (Metrics -> F18,
# today -> D6), leaving "today" as the active/tab-selected sheet.
$wsMetrics.Select()
$wsMetrics.Range("F18").Select()

$wsToday.Select()
$wsToday.Range("D6").Select()
